# Update computed profit/price figures across the Leviathan Profits workbook
# (sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1049.125
$ws.Range("I28").Value = 406.66666
$ws.Range("J28").Value = 2976.5
$ws.Range("K28").Value = 406.66666
$ws.Range("L28").Value = 2976.5
$ws.Range("M28").Value = 78.33334000000002
$ws.Range("N28").Value = -3946.5
$ws.Range("H40").Value = 4291.5835
$ws.Range("J40").Value = 4749.9
$ws.Range("L40").Value = 4749.9
$ws.Range("N40").Value = -5099.9
$ws.Range("H92").Value = 682.3
$ws.Range("I92").Value = 478.5
$ws.Range("J92").Value = 1497.5
$ws.Range("K92").Value = 478.5
$ws.Range("L92").Value = 1497.5
$ws.Range("M92").Value = 769.5
$ws.Range("N92").Value = -3993.5
$ws.Range("H113").Value = 5090.8184
$ws.Range("I113").Value = 3750
$ws.Range("J113").Value = 5857
$ws.Range("K113").Value = 3750
$ws.Range("L113").Value = 5857
$ws.Range("M113").Value = -496
$ws.Range("N113").Value = -12365
$ws.Range("H132").Value = 1485.8823
$ws.Range("I132").Value = 1284.1489
$ws.Range("J132").Value = 3856.25
$ws.Range("K132").Value = 3852.4467
$ws.Range("L132").Value = 11568.75
$ws.Range("M132").Value = -1322.4467
$ws.Range("N132").Value = -16628.75
$ws.Range("H134").Value = 131899.4
$ws.Range("J134").Value = 119999.664
$ws.Range("L134").Value = 119999.664
$ws.Range("N134").Value = -130139.664
$ws.Range("H137").Value = 2520.4167
$ws.Range("I137").Value = 2324.5
$ws.Range("J137").Value = 3500
$ws.Range("K137").Value = 6973.5
$ws.Range("L137").Value = 10500
$ws.Range("M137").Value = -4423.5
$ws.Range("N137").Value = -15600
$ws.Range("H138").Value = 2569.3215
$ws.Range("I138").Value = 2378.2856
$ws.Range("J138").Value = 3142.4285
$ws.Range("K138").Value = 7134.8568
$ws.Range("L138").Value = 9427.2855
$ws.Range("M138").Value = -1994.8568
$ws.Range("N138").Value = -19707.2855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H57").Value = 5178.5713
$ws.Range("I57").Value = 5178.5713
$ws.Range("K57").Value = 5178.5713
$ws.Range("M57").Value = -4694.5713
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("H126").Value = 5132.273
$ws.Range("I126").Value = 5132.273
$ws.Range("K126").Value = 15396.819
$ws.Range("M126").Value = -12926.819
$ws.Range("H132").Value = 1768.2258
$ws.Range("I132").Value = 1768.2258
$ws.Range("K132").Value = 5304.6774
$ws.Range("M132").Value = -2774.6774

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 80624
$ws.Range("J95").Value = 80624
$ws.Range("L95").Value = 80624
$ws.Range("N95").Value = -86116
$ws.Range("H113").Value = 5036.3
$ws.Range("I113").Value = 5036.3
$ws.Range("K113").Value = 5036.3
$ws.Range("M113").Value = -2866.3
$ws.Range("H134").Value = 715.7143
$ws.Range("I134").Value = 718.8333
$ws.Range("K134").Value = 2156.4999
$ws.Range("M134").Value = 378.5001000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 285761.84
$ws.Range("I2").Value = 53.25
$ws.Range("J2").Value = 666706.7
$ws.Range("K2").Value = 53.25
$ws.Range("L2").Value = 666706.7
$ws.Range("M2").Value = 59.75
$ws.Range("N2").Value = -666932.7
$ws.Range("H22").Value = 604.5
$ws.Range("I22").Value = 604.5
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 604.5
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -254.5
$ws.Range("H76").Value = 4919.5
$ws.Range("I76").Value = 4919.5
$ws.Range("K76").Value = 4919.5
$ws.Range("M76").Value = -4604.5
$ws.Range("H79").Value = 4919.5
$ws.Range("I79").Value = 4919.5
$ws.Range("K79").Value = 4919.5
$ws.Range("M79").Value = -3827.5
$ws.Range("H86").Value = 4887.5
$ws.Range("I86").Value = 4792.5
$ws.Range("K86").Value = 4792.5
$ws.Range("M86").Value = -3669.5
$ws.Range("H89").Value = 4887.5
$ws.Range("I89").Value = 4792.5
$ws.Range("K89").Value = 23962.5
$ws.Range("M89").Value = -18346.5
$ws.Range("H97").Value = 25098.5
$ws.Range("J97").Value = 25098.5
$ws.Range("L97").Value = 25098.5
$ws.Range("N97").Value = -27080.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 10538.947
$ws.Range("I56").Value = 10538.947
$ws.Range("K56").Value = 10538.947
$ws.Range("M56").Value = -10008.947
$ws.Range("H62").Value = 8799
$ws.Range("J62").Value = 9748.75
$ws.Range("L62").Value = 29246.25
$ws.Range("N62").Value = -30618.25
$ws.Range("H64").Value = 3044.2222
$ws.Range("I64").Value = 800
$ws.Range("J64").Value = 3324.75
$ws.Range("K64").Value = 2400
$ws.Range("L64").Value = 9974.25
$ws.Range("M64").Value = -2130
$ws.Range("N64").Value = -10514.25
$ws.Range("H65").Value = 8799
$ws.Range("J65").Value = 9748.75
$ws.Range("L65").Value = 87738.75
$ws.Range("N65").Value = -94602.75
$ws.Range("H67").Value = 3044.2222
$ws.Range("I67").Value = 800
$ws.Range("J67").Value = 3324.75
$ws.Range("K67").Value = 2400
$ws.Range("L67").Value = 9974.25
$ws.Range("M67").Value = -1464
$ws.Range("N67").Value = -11846.25
$ws.Range("H80").Value = 2400
$ws.Range("J80").Value = 2400
$ws.Range("L80").Value = 7200
$ws.Range("N80").Value = -9072
$ws.Range("H83").Value = 2400
$ws.Range("J83").Value = 2400
$ws.Range("L83").Value = 21600
$ws.Range("N83").Value = -30960
$ws.Range("H109").Value = 1499.5
$ws.Range("I109").Value = 1499.5
$ws.Range("K109").Value = 4498.5
$ws.Range("M109").Value = -3458.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3126.0232
$ws.Range("I132").Value = 1953.875
$ws.Range("K132").Value = 5861.625
$ws.Range("M132").Value = -3331.625
$ws.Range("H135").Value = 200000
$ws.Range("J135").Value = 200000
$ws.Range("L135").Value = 200000
$ws.Range("N135").Value = -210140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8193.675999999999
$ws.Range("I7").Value = 11197.823
$ws.Range("J7").Value = 5640.15
$ws.Range("K7").Value = 11197.823
$ws.Range("L7").Value = 5640.15
$ws.Range("M7").Value = -11085.823
$ws.Range("N7").Value = -5864.15
$ws.Range("H22").Value = 1413.4286
$ws.Range("I22").Value = 548.75
$ws.Range("J22").Value = 2566.3333
$ws.Range("K22").Value = 548.75
$ws.Range("L22").Value = 2566.3333
$ws.Range("M22").Value = -253.75
$ws.Range("N22").Value = -3156.3333
$ws.Range("H27").Value = 1413.4286
$ws.Range("I27").Value = 548.75
$ws.Range("J27").Value = 2566.3333
$ws.Range("K27").Value = 548.75
$ws.Range("L27").Value = 2566.3333
$ws.Range("M27").Value = -441.75
$ws.Range("N27").Value = -2780.3333
$ws.Range("H40").Value = 3915.805
$ws.Range("I40").Value = 3295.484
$ws.Range("J40").Value = 5838.8
$ws.Range("K40").Value = 3295.484
$ws.Range("L40").Value = 5838.8
$ws.Range("M40").Value = -3159.484
$ws.Range("N40").Value = -6110.8
$ws.Range("H105").Value = 36267.168
$ws.Range("J105").Value = 36267.168
$ws.Range("L105").Value = 36267.168
$ws.Range("N105").Value = -43255.168
$ws.Range("H126").Value = 8193.675999999999
$ws.Range("I126").Value = 11197.823
$ws.Range("J126").Value = 5640.15
$ws.Range("K126").Value = 33593.469
$ws.Range("L126").Value = 16920.45
$ws.Range("M126").Value = -31123.469
$ws.Range("N126").Value = -21860.45

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 9998.833000000001
$ws.Range("J25").Value = 9998.833000000001
$ws.Range("L25").Value = 9998.833000000001
$ws.Range("N25").Value = -10584.833
$ws.Range("H62").Value = 13042.5
$ws.Range("J62").Value = 12389.9
$ws.Range("L62").Value = 12389.9
$ws.Range("N62").Value = -13637.9
$ws.Range("H65").Value = 13042.5
$ws.Range("J65").Value = 12389.9
$ws.Range("L65").Value = 61949.5
$ws.Range("N65").Value = -68189.5
$ws.Range("H81").Value = 6550.143
$ws.Range("I81").Value = 918.6667
$ws.Range("K81").Value = 1837.3334
$ws.Range("M81").Value = -776.3334
$ws.Range("H84").Value = 6550.143
$ws.Range("I84").Value = 918.6667
$ws.Range("K84").Value = 9186.666999999999
$ws.Range("M84").Value = -3882.666999999999
$ws.Range("H107").Value = 20000808
$ws.Range("I107").Value = 599.65
$ws.Range("J107").Value = 100001650
$ws.Range("K107").Value = 1798.95
$ws.Range("L107").Value = 300004950
$ws.Range("M107").Value = 121.0500000000002
$ws.Range("N107").Value = -300008790
$ws.Range("H113").Value = 482.36365
$ws.Range("I113").Value = 424.94116
$ws.Range("K113").Value = 1274.82348
$ws.Range("M113").Value = 895.17652
$ws.Range("H126").Value = 1831.1
$ws.Range("I126").Value = 1648.92
$ws.Range("K126").Value = 4946.76
$ws.Range("M126").Value = -2476.76
$ws.Range("H132").Value = 4064.7437
$ws.Range("I132").Value = 4767.9
$ws.Range("K132").Value = 14303.7
$ws.Range("M132").Value = -11773.7

$ws2 = $wb.Worksheets.Item("ARM")
$ws2.Range("N108").ClearContents()
$ws2 = $wb.Worksheets.Item("CRP")
$ws2.Range("N22").ClearContents()
